# Kmeans grouper and labeler: update sentiment/user labels and append new tweet rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# Row 5: sentiment POSITIVE -> NEUTRAL
$ws.Range("E5").Value = "NEUTRAL"

# Row 6: user nas -> elonmusk
$ws.Range("D6").Value = "elonmusk"

# Row 7: user nas -> elonmusk, sentiment POSITIVE -> NEUTRAL
$ws.Range("D7").Value = "elonmusk"
$ws.Range("E7").Value = "NEUTRAL"

# Row 8: user nas -> elonmusk
$ws.Range("D8").Value = "elonmusk"

# --- Append new rows 9-11 ---
# Copy the formatting of the A column id cells (bold/border/center style) down
$ws.Range("A8").Copy()
$ws.Range("A9:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "twitter is amazing"
$ws.Range("C9").Value = "4,5,6"
$ws.Range("D9").Value = "nas"
$ws.Range("E9").Value = "POSITIVE"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "I hate twitter."
$ws.Range("C10").Value = "4,5,6"
$ws.Range("D10").Value = "nas"
$ws.Range("E10").Value = "NEGATIVE"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "I dunno"
$ws.Range("C11").Value = "4,5,6"
$ws.Range("D11").Value = "nas"
$ws.Range("E11").Value = "NEGATIVE"
